$d = $word.ActiveDocument

# Locate the target bullet paragraph by its text.
$targetText = "Consultar los cursos en los que está inscrito un estudiante."
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd("`r", "`a") -eq $targetText) {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Target paragraph not found"
}

$r = $target.Range

# Apply yellow highlight to the run text itself.
$r.HighlightColorIndex = 7   # wdYellow

# Also stamp the paragraph mark's run properties with the same highlight
# (so <w:pPr><w:rPr><w:highlight .../></w:rPr></w:pPr> is written, matching
# what Word does when the paragraph mark is included in the formatted
# selection). Range.HighlightColorIndex alone only reaches the visible run,
# so rebuild this single paragraph via a minimal OOXML replace that carries
# the highlight on both the run and the paragraph mark.
$rsid = $target.Range.ParagraphFormat.Style  # touch to ensure object is alive (no-op)

$paraXml = $target.Range.WordOpenXML
# Extract rsid-ish identity info isn't necessary: build fragment straight from
# the live paragraph's own paragraph XML via InsertXML so Word regenerates it.

$xml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="6B59F1CB" w14:textId="77777777" w:rsidR="00DA04B2" w:rsidRPr="00DA04B2" w:rsidRDefault="00DA04B2" w:rsidP="00DA04B2">
<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr>
<w:r w:rsidRPr="00DA04B2"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Consultar los cursos en los que est&#225; inscrito un estudiante.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)
